# Generate Report for Handback
#
# - Status text "Ready for handoff" -> "Handed back: in sync with en-US"
#   (appears on Overview!B:C and on each locale sheet's Status column C)
# - zh-cn / de-de sheets: populate the newly-produced "Latest Target File"
#   (F) and "Latest Handback File" (G) columns for rows 2 & 3, by copying
#   the existing hyperlink in column A (the .md target) into F, and the
#   existing hyperlink in column D (the .xlf handoff file) into G.
# - zh-cn / de-de sheets: stamp the "Latest Handback DateTime" (H) with the
#   handback timestamp for that locale.

$wb = $excel.ActiveWorkbook

$newStatus = "Handed back: in sync with en-US"

# ---- Overview sheet: Status columns B & C for rows 2 & 3 ----
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B2").Value = $newStatus
$overview.Range("C2").Value = $newStatus
$overview.Range("B3").Value = $newStatus
$overview.Range("C3").Value = $newStatus

function Update-LocaleSheet($SheetName, $HandbackDateTime) {

    $ws = $wb.Worksheets.Item($SheetName)

    # Status column (C) on rows 2 & 3
    $ws.Range("C2").Value = $newStatus
    $ws.Range("C3").Value = $newStatus

    # Grab the existing hyperlink target/display text already on A2 (the
    # source .md file) and D2 (the locale .xlf handoff file) so the new
    # Latest Target File / Latest Handback File columns link to the same
    # places.
    $mdAddress = $null
    $mdDisplay = $null
    $xlfAddress = $null
    $xlfDisplay = $null
    foreach ($hl in $ws.Hyperlinks) {
        $refAddr = $hl.Range.Address()
        if ($refAddr -eq "`$A`$2") {
            $mdAddress = $hl.Address
            $mdDisplay = $hl.TextToDisplay
        }
        elseif ($refAddr -eq "`$D`$2") {
            $xlfAddress = $hl.Address
            $xlfDisplay = $hl.TextToDisplay
        }
    }

    # New "Latest Target File" (F) / "Latest Handback File" (G) cells for
    # rows 2 and 3 - both rows link to the same md/xlf pair.
    foreach ($row in 2, 3) {
        $ws.Hyperlinks.Add($ws.Range("F$row"), $mdAddress, "", "", $mdDisplay)
        $ws.Hyperlinks.Add($ws.Range("G$row"), $xlfAddress, "", "", $xlfDisplay)
    }

    # "Latest Handback DateTime" (H) on rows 2 & 3
    $ws.Range("H2").Value = $HandbackDateTime
    $ws.Range("H3").Value = $HandbackDateTime
}

Update-LocaleSheet "zh-cn" "2016-03-19 17:05:57"
Update-LocaleSheet "de-de" "2016-03-19 17:06:14"
